$d = $word.ActiveDocument
$cr = [string][char]13

$tinderbox = "Tinderbox Marine Reserve is a small no-take area founded in 1991 near Tasmania" + [string][char]8217 + "s capital city of Hobart.  Despite its small size, the reserve has affected many species in interesting ways, including a ten-fold increase in the number of larger fish (>300mm) inside the reserve compared to the fished areas outside.  The population of large bastard trumpeter fish benefitted the most, increasing by a factor of eight inside the reserve over the course of the study.  However, the population of smaller trumpeters remained the same, a result indicative of high fishing pressure outside the reserve that removed the larger individuals of the population.  The reserve had similar effects on the southern rock lobster, allowing lobsters in the reserve to grow much larger while increasing in density.  Additional information from "

$fernando = "This island archipelago 345km northeast of Brazil is a very important nursery area for juvenile Caribbean reef sharks.  Fortunately for the protection efforts of these valuable members of the tropical ecosystem, Fernando de Noronha is surrounded by a marine reserve covering 70% of its coastal waters.  A tracking study showed that juvenile reef sharks spend most of their time in the reserve area, rarely venturing out of the reserve boundaries into the unprotected area where human impacts like fishing pressure and boating noise are highest.  In fact, not a single shark during the duration of the study crossed from one side of the unprotected area to the other, possibly indicating an important preference for undisturbed habitat during one of their most vulnerable life stages.  While adult sharks roam across much larger areas, Fernando de Noronha has shown that marine reserves can help protect some shark populations by providing a sanctuary for the individuals that will become the next generation of adult sharks"

# Insert right before the trailing "." of "died back." (i.e. two positions
# before the very end of the document body, since position End-1 is right
# after the "." and before the final paragraph mark). This way the existing
# "." stays as the very last character of the document, now terminating the
# new final sentence ("...adult sharks.") instead of "died back".
$insertPoint = $d.Content.End - 2
$r = $d.Range($insertPoint, $insertPoint)

$full = "." + $cr + $cr + "Tinderbox Marine Reserve, Tasmania" + $cr + $cr + $tinderbox + $cr + $cr + "Fernando de Noronha, Brazil" + $cr + $cr + $fernando

$r.InsertAfter($full)
